# Auto-generated edit script: applies updated market-price snapshot values
# to the leve-profit tracker sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 11153.777
$ws.Range("I76").Value = 12925.714
$ws.Range("K76").Value = 12925.714
$ws.Range("M76").Value = -12610.714
$ws.Range("H79").Value = 11153.777
$ws.Range("I79").Value = 12925.714
$ws.Range("K79").Value = 12925.714
$ws.Range("M79").Value = -11833.714
$ws.Range("H131").Value = 2838.7856
$ws.Range("I131").Value = 2337.7778
$ws.Range("J131").Value = 3740.6
$ws.Range("K131").Value = 7013.3334
$ws.Range("L131").Value = 11221.8
$ws.Range("M131").Value = -1973.3334
$ws.Range("N131").Value = -21301.8
$ws.Range("H132").Value = 1874.3422
$ws.Range("I132").Value = 1603.625
$ws.Range("J132").Value = 3318.1667
$ws.Range("K132").Value = 4810.875
$ws.Range("L132").Value = 9954.500100000001
$ws.Range("M132").Value = -2280.875
$ws.Range("N132").Value = -15014.5001
$ws.Range("H137").Value = 3094.7827
$ws.Range("I137").Value = 3549.375
$ws.Range("K137").Value = 10648.125
$ws.Range("M137").Value = -8098.125
$ws.Range("H141").Value = 2884.875
$ws.Range("I141").Value = 2216.8
$ws.Range("K141").Value = 6650.400000000001
$ws.Range("M141").Value = -1470.400000000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2723614
$ws.Range("I32").Value = 2846107.5
$ws.Range("K32").Value = 2846107.5
$ws.Range("M32").Value = -2845820.5
$ws.Range("H45").Value = 5107.778
$ws.Range("I45").Value = 1660.25
$ws.Range("K45").Value = 1660.25
$ws.Range("M45").Value = -1283.25
$ws.Range("H61").Value = 7245.6313
$ws.Range("I61").Value = 3419.4
$ws.Range("K61").Value = 3419.4
$ws.Range("M61").Value = -3207.4
$ws.Range("H74").Value = 51350.332
$ws.Range("I74").Value = 68672.086
$ws.Range("K74").Value = 68672.086
$ws.Range("M74").Value = -67798.086
$ws.Range("H77").Value = 51350.332
$ws.Range("I77").Value = 68672.086
$ws.Range("K77").Value = 343360.43
$ws.Range("M77").Value = -338992.43
$ws.Range("H122").Value = 15200.941
$ws.Range("I122").Value = 21436.9
$ws.Range("K122").Value = 64310.7
$ws.Range("M122").Value = -61860.7
$ws.Range("H132").Value = 3771.5745
$ws.Range("I132").Value = 1524.7576
$ws.Range("K132").Value = 4574.2728
$ws.Range("M132").Value = -2044.2728
$ws.Range("H136").Value = 7245.6313
$ws.Range("I136").Value = 3419.4
$ws.Range("K136").Value = 10258.2
$ws.Range("M136").Value = -7708.200000000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 30490642
$ws.Range("I86").Value = 14708378
$ws.Range("J86").Value = 41669744
$ws.Range("K86").Value = 14708378
$ws.Range("L86").Value = 41669744
$ws.Range("M86").Value = -14707255
$ws.Range("N86").Value = -41671990
$ws.Range("H89").Value = 30490642
$ws.Range("I89").Value = 14708378
$ws.Range("J89").Value = 41669744
$ws.Range("K89").Value = 73541890
$ws.Range("L89").Value = 208348720
$ws.Range("M89").Value = -73536274
$ws.Range("N89").Value = -208359952
$ws.Range("H94").Value = 1571.5186
$ws.Range("I94").Value = 714.1818
$ws.Range("J94").Value = 5343.8
$ws.Range("K94").Value = 714.1818
$ws.Range("L94").Value = 5343.8
$ws.Range("M94").Value = -263.1818
$ws.Range("N94").Value = -6245.8
$ws.Range("H99").Value = 3249359.2
$ws.Range("J99").Value = 9094202
$ws.Range("L99").Value = 9094202
$ws.Range("N99").Value = -9097198
$ws.Range("H107").Value = 40181624
$ws.Range("I107").Value = 51138480
$ws.Range("J107").Value = 6487.1665
$ws.Range("K107").Value = 51138480
$ws.Range("L107").Value = 6487.1665
$ws.Range("M107").Value = -51136560
$ws.Range("N107").Value = -10327.1665
$ws.Range("H134").Value = 5535.1777
$ws.Range("I134").Value = 2100.862
$ws.Range("K134").Value = 6302.586
$ws.Range("M134").Value = -3767.586
$ws.Range("H139").Value = 29599.8
$ws.Range("J139").Value = 29599.8
$ws.Range("L139").Value = 29599.8
$ws.Range("N139").Value = -39879.8
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11209.29
$ws.Range("I31").Value = 5097.385
$ws.Range("J31").Value = 15623.444
$ws.Range("K31").Value = 5097.385
$ws.Range("L31").Value = 15623.444
$ws.Range("M31").Value = -4802.385
$ws.Range("N31").Value = -16213.444
$ws.Range("H34").Value = 11209.29
$ws.Range("I34").Value = 5097.385
$ws.Range("J34").Value = 15623.444
$ws.Range("K34").Value = 5097.385
$ws.Range("L34").Value = 15623.444
$ws.Range("M34").Value = -4895.385
$ws.Range("N34").Value = -16027.444
$ws.Range("H60").Value = 19998.334
$ws.Range("I60").Value = 10000
$ws.Range("K60").Value = 10000
$ws.Range("M60").Value = -9489
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("M123").ClearContents()
$ws.Range("H132").Value = 5899
$ws.Range("I132").Value = 1839.875
$ws.Range("K132").Value = 5519.625
$ws.Range("M132").Value = -2989.625
$ws.Range("H135").Value = 88000
$ws.Range("J135").Value = 88000
$ws.Range("L135").Value = 88000
$ws.Range("N135").Value = -98140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 134703.86
$ws.Range("I2").Value = 50.272728
$ws.Range("K2").Value = 301.636368
$ws.Range("M2").Value = -188.636368
$ws.Range("H48").Value = 9000
$ws.Range("J48").Value = 9000
$ws.Range("L48").Value = 27000
$ws.Range("N48").Value = -27500
$ws.Range("H117").Value = 1000
$ws.Range("I117").Value = 0
$ws.Range("K117").Value = 0
$ws.Range("M117").ClearContents()
$ws.Range("H136").Value = 1850
$ws.Range("I136").Value = 1850
$ws.Range("K136").Value = 5550
$ws.Range("M136").Value = -450
$ws.Range("H137").Value = 401394.6
$ws.Range("I137").Value = 251743.5
$ws.Range("K137").Value = 755230.5
$ws.Range("M137").Value = -750130.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10169.25
$ws.Range("I70").Value = 9270.5
$ws.Range("K70").Value = 9270.5
$ws.Range("M70").Value = -9000.5
$ws.Range("H73").Value = 10169.25
$ws.Range("I73").Value = 9270.5
$ws.Range("K73").Value = 9270.5
$ws.Range("M73").Value = -8334.5
$ws.Range("H107").Value = 768.1539
$ws.Range("I107").Value = 578.9
$ws.Range("K107").Value = 578.9
$ws.Range("M107").Value = 1341.1
$ws.Range("H132").Value = 3762.795
$ws.Range("I132").Value = 1490.84
$ws.Range("K132").Value = 4472.52
$ws.Range("M132").Value = -1942.52
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2833.125
$ws.Range("I22").Value = 1500
$ws.Range("J22").Value = 3439.0908
$ws.Range("K22").Value = 1500
$ws.Range("L22").Value = 3439.0908
$ws.Range("M22").Value = -1205
$ws.Range("N22").Value = -4029.0908
$ws.Range("H27").Value = 2833.125
$ws.Range("I27").Value = 1500
$ws.Range("J27").Value = 3439.0908
$ws.Range("K27").Value = 1500
$ws.Range("L27").Value = 3439.0908
$ws.Range("M27").Value = -1393
$ws.Range("N27").Value = -3653.0908
$ws.Range("H46").Value = 1232.3043
$ws.Range("I46").Value = 272.5
$ws.Range("J46").Value = 1434.3684
$ws.Range("K46").Value = 272.5
$ws.Range("L46").Value = 1434.3684
$ws.Range("M46").Value = -84.5
$ws.Range("N46").Value = -1810.3684
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("H133").Value = 35833.25
$ws.Range("J133").Value = 35833.25
$ws.Range("L133").Value = 35833.25
$ws.Range("N133").Value = -40893.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 24996
$ws.Range("I15").Value = 24996
$ws.Range("K15").Value = 24996
$ws.Range("M15").Value = -24708
$ws.Range("H124").Value = 24428
$ws.Range("J124").Value = 24428
$ws.Range("L124").Value = 24428
$ws.Range("N124").Value = -34248
$ws.Range("H132").Value = 31274882
$ws.Range("I132").Value = 55567710
$ws.Range("J132").Value = 41242.715
$ws.Range("K132").Value = 166703130
$ws.Range("L132").Value = 123728.145
$ws.Range("M132").Value = -166700600
$ws.Range("N132").Value = -128788.145
$ws.Range("H136").Value = 43485904
$ws.Range("I136").Value = 166668930
$ws.Range("J136").Value = 9538.471
$ws.Range("K136").Value = 500006790
$ws.Range("L136").Value = 28615.413
$ws.Range("M136").Value = -500004240
$ws.Range("N136").Value = -33715.413
